$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(866, 8).Value = 503
$ws.Cells.Item(867, 8).Value = 536
$ws.Cells.Item(868, 8).Value = 548
$ws.Cells.Item(869, 8).Value = 567
$ws.Cells.Item(870, 8).Value = 593
$ws.Cells.Item(871, 8).Value = 644
$ws.Cells.Item(872, 8).Value = 700
$ws.Cells.Item(873, 8).Value = 724
$ws.Cells.Item(874, 8).Value = 745
$ws.Cells.Item(875, 8).Value = 762
$ws.Cells.Item(876, 8).Value = 741
$ws.Cells.Item(877, 8).Value = 717
$ws.Cells.Item(878, 8).Value = 725
$ws.Cells.Item(879, 8).Value = 753
$ws.Cells.Item(880, 8).Value = 725
$ws.Cells.Item(881, 8).Value = 702
$ws.Cells.Item(882, 8).Value = 659
$ws.Cells.Item(883, 8).Value = 663
$ws.Cells.Item(884, 8).Value = 590
$ws.Cells.Item(885, 8).Value = 616
$ws.Cells.Item(886, 8).Value = 654
$ws.Cells.Item(887, 8).Value = 603
$ws.Cells.Item(888, 8).Value = 566
$ws.Cells.Item(889, 8).Value = 551
$ws.Cells.Item(890, 8).Value = 532
$ws.Cells.Item(891, 8).Value = 493
$ws.Cells.Item(892, 8).Value = 517
$ws.Cells.Item(893, 8).Value = 543
$ws.Cells.Item(894, 8).Value = 517
$ws.Cells.Item(895, 8).Value = 513
$ws.Cells.Item(896, 8).Value = 513
$ws.Cells.Item(897, 8).Value = 525
$ws.Cells.Item(898, 8).Value = 499
$ws.Cells.Item(899, 8).Value = 521
$ws.Cells.Item(900, 8).Value = 555
$ws.Cells.Item(901, 8).Value = 554
$ws.Cells.Item(902, 8).Value = 516
$ws.Cells.Item(903, 8).Value = 473
$ws.Cells.Item(904, 6).Value = 5380
$ws.Cells.Item(904, 8).Value = 479
$ws.Cells.Item(905, 8).Value = 442
$ws.Cells.Item(906, 8).Value = 444
$ws.Cells.Item(907, 8).Value = 486
$ws.Cells.Item(908, 8).Value = 503
$ws.Cells.Item(909, 8).Value = 456
$ws.Cells.Item(910, 6).Value = 2794
$ws.Cells.Item(910, 8).Value = 411
$ws.Cells.Item(911, 8).Value = 420
$ws.Cells.Item(912, 6).Value = 2959
$ws.Cells.Item(912, 8).Value = 382
$ws.Cells.Item(913, 8).Value = 401
$ws.Cells.Item(914, 8).Value = 419
$ws.Cells.Item(915, 8).Value = 403
$ws.Cells.Item(916, 8).Value = 394
$ws.Cells.Item(917, 8).Value = 386
$ws.Cells.Item(918, 6).Value = 3620
$ws.Cells.Item(918, 7).Value = 247
$ws.Cells.Item(918, 8).Value = 402
$ws.Cells.Item(919, 6).Value = 3894
$ws.Cells.Item(919, 8).Value = 374
$ws.Cells.Item(920, 8).Value = 385
$ws.Cells.Item(921, 6).Value = 1016
$ws.Cells.Item(921, 7).Value = 103
$ws.Cells.Item(921, 8).Value = 396
$ws.Cells.Item(922, 6).Value = 6055
$ws.Cells.Item(922, 8).Value = 368
$ws.Cells.Item(923, 8).Value = 366
$ws.Cells.Item(924, 6).Value = 2654
$ws.Cells.Item(924, 8).Value = 358
$ws.Cells.Item(925, 6).Value = 1171
$ws.Cells.Item(925, 8).Value = 376
$ws.Cells.Item(926, 6).Value = 2970
$ws.Cells.Item(926, 8).Value = 366
$ws.Cells.Item(927, 8).Value = 396
$ws.Cells.Item(928, 8).Value = 426
$ws.Cells.Item(929, 6).Value = 5005
$ws.Cells.Item(929, 7).Value = 567
$ws.Cells.Item(929, 8).Value = 434
$ws.Cells.Item(930, 6).Value = 3528
$ws.Cells.Item(930, 7).Value = 395
$ws.Cells.Item(930, 8).Value = 428
$ws.Cells.Item(931, 6).Value = 3789
$ws.Cells.Item(931, 7).Value = 321
$ws.Cells.Item(931, 8).Value = 437
$ws.Cells.Item(932, 6).Value = 3349
$ws.Cells.Item(932, 7).Value = 309
$ws.Cells.Item(932, 8).Value = 447
$ws.Cells.Item(933, 6).Value = 2517
$ws.Cells.Item(933, 7).Value = 316
$ws.Cells.Item(933, 8).Value = 443
$ws.Cells.Item(934, 6).Value = 642
$ws.Cells.Item(934, 7).Value = 95
$ws.Cells.Item(934, 8).Value = 443
$ws.Cells.Item(935, 6).Value = 350
$ws.Cells.Item(935, 7).Value = 61
$ws.Cells.Item(935, 8).Value = 443
